$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SubjectSampleMapping_DD")

# Update the VARDESC column (column B) for the two data rows with the
# full descriptive text (replacing the short "Subject ID"/"Sample ID" labels).
$ws.Range("B2").Value = "A sequence of letters, numbers, or characters that uniquely identifies the subject who has taken part in the investigation or research study."
$ws.Range("B3").Value = "A unique sequence of alphanumeric characters used to identify the specimen at it's point of origin."

# Widen columns B and C to fit the new, longer text and drop the old
# "best fit" auto-sizing in favor of an explicit, equal custom width
# (~103.33 characters, matching the author's manual resize).
$ws.Range("B:C").ColumnWidth = 102.42

# Clear the selection stored in the sheet view so it no longer points at A3
# (Excel omits the <selection> element entirely once it is back at the
# default top-left cell A1).
$ws.Range("A1").Select() | Out-Null
